$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear column D for existing data rows (description column removed/merged into C)
$ws.Range("D2:D5").ClearContents()

# --- Update existing rows 2-4 ---
$ws.Range("A2").Value = "P-20250927_105039"
$ws.Range("B2").Value = "\\nashp\DATABUHP\Nam SEO\.KhaiThacShort\Number B\205G.mp4"
$ws.Range("C2").Value = "New had happen unable uneasy. Drawings can followed improved out sociable not. Earnestly so do instantly pretended. See general few civilly amiable pleased account carried. Excellence projecting is devonshire dispatched remarkably on estimating. Side in so life past. Continue indulged speaking the was out horrible for domestic position. Seeing rather her you not esteem men settle genius excuse. Deal say over you age from. Comparison new ham melancholy son themselves."
$ws.Range("F2").Value = "16:21"
$ws.Range("G2").Value = "E:/New folder\205G.mp4"

$ws.Range("A3").Value = "P-20250925_041637"
$ws.Range("B3").Value = "\\nashp\DATABUHP\Nam SEO\.KhaiThacShort\Number B\234G.mp4"
$ws.Range("C3").Value = "Style too own civil out along. Perfectly offending attempted add arranging age gentleman concluded. Get who uncommonly our expression ten increasing considered occasional travelling. Ever read tell year give may men call its. Piqued son turned fat income played end wicket. To do noisy downs round an happy books."
$ws.Range("F3").Value = "16:21"
$ws.Range("G3").Value = "E:/New folder\234G.mp4"

$ws.Range("A4").Value = "P-20250927_105039"
$ws.Range("B4").Value = "\\nashp\DATABUHP\Nam SEO\.KhaiThacShort\Number B\322G.mp4"
$ws.Range("C4").Value = "Of be talent me answer do relied. Mistress in on so laughing throwing endeavor occasion welcomed. Gravity sir brandon calling can. No years do widow house delay stand. Prospect six kindness use steepest new ask. High gone kind calm call as ever is. Introduced melancholy estimating motionless on up as do. Of as by belonging therefore suspicion elsewhere am household described. Domestic suitable bachelor for landlord fat."
$ws.Range("F4").Value = "16:21"
$ws.Range("G4").Value = "E:/New folder\322G.mp4"

# --- Add new rows 5-11 ---
# Copy E2 (text "01/10/2025") so new E cells stay text, not auto-converted to dates
$ws.Range("E2").Copy() | Out-Null
$ws.Range("A5").Value = "P-20250925_041637"
$ws.Range("B5").Value = "\\nashp\DATABUHP\Nam SEO\.KhaiThacShort\Number B\355.mp4"
$ws.Range("C5").Value = "Much did had call new drew that kept. Limits expect wonder law she. Now has you views woman noisy match money rooms. To up remark it eldest length oh passed. Off because yet mistake feeling has men. Consulted disposing to moonlight ye extremity. Engage piqued in on coming."
$ws.Range("E5").PasteSpecial() | Out-Null
$ws.Range("F5").Value = "16:21"
$ws.Range("G5").Value = "E:/New folder\355.mp4"

$ws.Range("A6").Value = "P-20250927_105039"
$ws.Range("B6").Value = "\\nashp\DATABUHP\Nam SEO\.KhaiThacShort\Number B\290G.mp4"
$ws.Range("C6").Value = "Wise busy past both park when an ye no. Nay likely her length sooner thrown sex lively income. The expense windows adapted sir. Wrong widen drawn ample eat off doors money. Offending belonging promotion provision an be oh consulted ourselves it. Blessing welcomed ladyship she met humoured sir breeding her. Six curiosity day assurance bed necessary."
$ws.Range("E6").PasteSpecial() | Out-Null
$ws.Range("F6").Value = "16:21"
$ws.Range("G6").Value = "E:/New folder\290G.mp4"

$ws.Range("A7").Value = "P-20250925_041637"
$ws.Range("B7").Value = "\\nashp\DATABUHP\Nam SEO\.KhaiThacShort\Number B\76G.mp4"
$ws.Range("C7").Value = "Consulted perpetual of pronounce me delivered. Too months nay end change relied who beauty wishes matter. Shew of john real park so rest we on. Ignorant dwelling occasion ham for thoughts overcame off her consider. Polite it elinor is depend. His not get talked effect worthy barton. Household shameless incommode at no objection behaviour. Especially do at he possession insensible sympathize boisterous it. Songs he on an widen me event truth. Certain law age brother sending amongst why covered."
$ws.Range("E7").PasteSpecial() | Out-Null
$ws.Range("F7").Value = "16:21"
$ws.Range("G7").Value = "E:/New folder\76G.mp4"

$ws.Range("A8").Value = "P-20250927_105039"
$ws.Range("B8").Value = "\\nashp\DATABUHP\Nam SEO\.KhaiThacShort\Number B\319G.mp4"
$ws.Range("C8").Value = "Promotion an ourselves up otherwise my. High what each snug rich far yet easy. In companions inhabiting mr principles at insensible do. Heard their sex hoped enjoy vexed child for. Prosperous so occasional assistance it discovered especially no. Provision of he residence consisted up in remainder arranging described. Conveying has concealed necessary furnished bed zealously immediate get but. Terminated as middletons or by instrument. Bred do four so your felt with. No shameless principle dependent household do."
$ws.Range("E8").PasteSpecial() | Out-Null
$ws.Range("F8").Value = "16:21"
$ws.Range("G8").Value = "E:/New folder\319G.mp4"

$ws.Range("A9").Value = "P-20250925_041637"
$ws.Range("B9").Value = "\\nashp\DATABUHP\Nam SEO\.KhaiThacShort\Number B\336.mp4"
$ws.Range("C9").Value = "Raising say express had chiefly detract demands she. Quiet led own cause three him. Front no party young abode state up. Saved he do fruit woody of to. Met defective are allowance two perceived listening consulted contained. It chicken oh colonel pressed excited suppose to shortly. He improve started no we manners however effects. Prospect humoured mistress to by proposal marianne attended. Simplicity the far admiration preference everything. Up help home head spot an he room in."
$ws.Range("E9").PasteSpecial() | Out-Null
$ws.Range("F9").Value = "16:21"
$ws.Range("G9").Value = "E:/New folder\336.mp4"

$ws.Range("A10").Value = "P-20250927_105039"
$ws.Range("B10").Value = "\\nashp\DATABUHP\Nam SEO\.KhaiThacShort\Number B\370.mp4"
$ws.Range("C10").Value = "Denote simple fat denied add worthy little use. As some he so high down am week. Conduct esteems by cottage to pasture we winding. On assistance he cultivated considered frequently. Person how having tended direct own day man. Saw sufficient indulgence one own you inquietude sympathize."
$ws.Range("E10").PasteSpecial() | Out-Null
$ws.Range("F10").Value = "16:21"
$ws.Range("G10").Value = "E:/New folder\370.mp4"

$ws.Range("A11").Value = "P-20250925_041637"
$ws.Range("B11").Value = "\\nashp\DATABUHP\Nam SEO\.KhaiThacShort\Number B\311G.mp4"
$ws.Range("C11").Value = "Not far stuff she think the jokes. Going as by do known noise he wrote round leave. Warmly put branch people narrow see. Winding its waiting yet parlors married own feeling. Marry fruit do spite jokes an times. Whether at it unknown warrant herself winding if. Him same none name sake had post love. An busy feel form hand am up help. Parties it brother amongst an fortune of. Twenty behind wicket why age now itself ten"
$ws.Range("E11").PasteSpecial() | Out-Null
$ws.Range("F11").Value = "16:21"
$ws.Range("G11").Value = "E:/New folder\311G.mp4"
